# report commesse e collaboratori - miglioramento interfaccia
# Add two summary formulas on row 1/2 column E, and move the selection to E2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E1: total cost summed from the detail rows (col F, rows 12-500)
$ws.Range("E1").Formula = "=SUM(F12:F500)"

# E2: total hours summed from the detail rows (col C, rows 12-500)
$ws.Range("E2").Formula = "=SUM(C12:C500)"

# Move/record the active selection on E2 (matches the updated sheetView)
[void]$ws.Range("E2").Select()
